# Adds Android/iOS related test-step columns (quantity tracking) to Sheet1
# and introduces a new "Questions" worksheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# Sheet1: insert quantity columns next to each item column, and extend the
# row-5 (Login test) data with list/item/quantity sample values plus a new
# Comment note. Also records the Pass result + comment for the register step.
# ---------------------------------------------------------------------------

# Header row (row 1): item1 stays in G1; quantity columns are interleaved
# between each item column, and Result/Comment shift two columns to the
# right (from J/K to M/N).
$ws1.Range("H1").Value = "quantity1"
$ws1.Range("I1").Value = "item2"
$ws1.Range("J1").Value = "quantity2"
$ws1.Range("K1").Value = "item3"
$ws1.Range("L1").Value = "quantity3"
$ws1.Range("M1").Value = "Result"
$ws1.Range("N1").Value = "Comment"

# Row 4 ("Register a user" step): record the expected result + comment.
$ws1.Range("M4").Value = "Pass"
$ws1.Range("N4").Value = "UI should display a message saying the user account is successfully created."

# Row 5 ("Login to user account" step): fill in sample list/item/quantity
# values exercised by the test, and a comment describing quantity behavior.
$ws1.Range("E5").Value = "password1001"
$ws1.Range("F5").Value = "listname1001"
$ws1.Range("G5").Value = "item11001"
$ws1.Range("H5").Value = 1
$ws1.Range("I5").Value = "item21001"
$ws1.Range("J5").Value = 2
$ws1.Range("K5").Value = "item31001"
$ws1.Range("L5").Value = 3
$ws1.Range("N5").Value = "Quantity should be defaulted to 1. `nQuantity should be increased or decreased with + and - buttons"
$ws1.Range("N5").WrapText = $true
$ws1.Rows("5").RowHeight = 29

# Column widths: make room for the new quantity columns and the wider
# Comment column.
$ws1.Columns("F").ColumnWidth = 11.3317
$ws1.Columns("G").ColumnWidth = 8.6649
$ws1.Range("H1:J1").ColumnWidth = 8.6649
$ws1.Columns("K").ColumnWidth = 8.6649
$ws1.Columns("L").ColumnWidth = 4.667
$ws1.Columns("N").ColumnWidth = 66.0002

$ws1.Range("K5").Select() | Out-Null

# ---------------------------------------------------------------------------
# New "Questions" worksheet, placed after Sheet1.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Questions"

$ws2.Range("A1").Value = "#"
$ws2.Range("B1").Value = "Question"
$ws2.Range("A2").Value = 1
$ws2.Range("B2").Value = "How to make the app available in all the regions of the world?"

$ws2.Columns("B").ColumnWidth = 52.667

$ws2.Range("B15").Select() | Out-Null

# Leave Sheet1 as the active/visible sheet when the workbook is reopened.
$ws1.Activate() | Out-Null
$ws1.Range("K5").Select() | Out-Null
